$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.579.42"
$ws.Range("E2").Value = "  +2.27%  "
$ws.Range("D3").Value = "3.542.37"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'609.81"
$ws.Range("E5").Value = "  +5.17%  "
$ws.Range("D6").Value = "'173.30"
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("E7").Value = "  +0.89%  "
$ws.Range("D8").Value = "3.538.70"
$ws.Range("E8").Value = "  +1.49%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  +6.43%  "
$ws.Range("D11").Value = "'6.74"
$ws.Range("E11").Value = "  +1.35%  "
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("D13").Value = "'47.38"
$ws.Range("E13").Value = "  +2.02%  "
$ws.Range("D14").Value = "'0.0000280"
$ws.Range("E14").Value = "  +2.25%  "
$ws.Range("D15").Value = "4.105.98"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").Value = "'627.47"
$ws.Range("E16").Value = "  -6.57%  "
$ws.Range("E17").Value = "  -2.67%  "
$ws.Range("D18").Value = "70.573.87"
$ws.Range("E18").Value = "  +2.37%  "
$ws.Range("D19").Value = "3.541.99"
$ws.Range("E19").Value = "  +1.49%  "
$ws.Range("E20").Value = "  -1.63%  "
$ws.Range("D21").Value = "'17.40"
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'9.99"
$ws.Range("E22").Value = "  -10.07%  "
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").Value = "'0.887"
$ws.Range("E23").Value = "  -1.02%  "
$ws.Range("D24").Value = "'15.91"
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("D25").Value = "'96.82"
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("D26").Value = "'3.86"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").Value = "'2.62"
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("D30").Value = "'33.47"
$ws.Range("E30").Value = "  +2.29%  "
$ws.Range("D31").Value = "'8.49"
$ws.Range("E31").Value = "  -1.65%  "
$ws.Range("E32").Value = "  -2.28%  "
$ws.Range("D33").Value = "'1.33"
$ws.Range("E33").Value = "  -1.45%  "
$ws.Range("D34").Value = "'7.00"
$ws.Range("E34").Value = "  -2.90%  "
$ws.Range("D35").Value = "'569.57"
$ws.Range("E35").Value = "  -3.75%  "
$ws.Range("D36").Value = "'3.62"
$ws.Range("E36").Value = "  +1.74%  "
$ws.Range("D37").Value = "'10.80"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").Value = "'57.71"
$ws.Range("E38").Value = "  +1.45%  "
$ws.Range("E39").Value = "  -1.27%  "
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("D41").Value = "'0.0459"
$ws.Range("E41").Value = "  +5.56%  "
$ws.Range("D43").Value = "'0.329"
$ws.Range("E43").Value = "  -1.23%  "
$ws.Range("D44").Value = "3.348.11"
$ws.Range("E44").Value = "  -1.47%  "
$ws.Range("E45").Value = "  +5.59%  "
$ws.Range("E46").Value = "  +1.76%  "
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("E48").Value = "  +2.90%  "
$ws.Range("D49").Value = "'0.130"
$ws.Range("E49").Value = "  -1.62%  "
$ws.Range("D50").Value = "'133.72"
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("E51").Value = "  +1.03%  "
